$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.109.92"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "3.052.63"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'585.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").Value = "'151.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.18%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.536"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "3.052.29"
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("D11").Value = "'5.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("E12").Value = "  -2.73%  "
$ws.Range("E13").Value = "  -2.44%  "
$ws.Range("D14").Value = "'36.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.51%  "
$ws.Range("E15").Value = "  +1.87%  "
$ws.Range("D16").Value = "3.555.57"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").Value = "63.091.04"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").Value = "3.055.06"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").Value = "'477.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("D21").Value = "'14.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.51%  "
$ws.Range("E22").Value = "  -1.47%  "
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("E24").Value = "  +1.29%  "
$ws.Range("D25").Value = "'82.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.74%  "
$ws.Range("D26").Value = "'12.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.96%  "
$ws.Range("E27").Value = "  +6.10%  "
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("E29").Value = "  +1.25%  "
$ws.Range("D30").Value = "'2.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").Value = "'27.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.74%  "
$ws.Range("E34").Value = "  -2.60%  "
$ws.Range("E35").Value = "  +1.17%  "
$ws.Range("D36").Value = "0.0₃0816"
$ws.Range("E36").Value = "  -2.88%  "
$ws.Range("D37").Value = "'3.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.36%  "
$ws.Range("E38").Value = "  -3.33%  "
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("D40").Value = "'9.24"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").Value = "'50.42"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D42").Value = "'433.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.18%  "
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("E44").Value = "  +2.87%  "
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("D46").Value = "2.827.51"
$ws.Range("E46").Value = "  +1.24%  "
$ws.Range("D47").Value = "'38.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.22%  "
$ws.Range("D48").Value = "'129.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").Value = "'25.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("E51").Value = "  -1.59%  "
